$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 4,5,6 hold F values that changed
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 1537
$wsExhibition.Range("F5").Value = 711
$wsExhibition.Range("F6").Value = 20

# Sheet "全部类型" (All types) - rows 4,6,7 hold F values that changed
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 1537
$wsAllTypes.Range("F6").Value = 711
$wsAllTypes.Range("F7").Value = 20
